$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (copying row 2's formatting/types down),
# shifting existing rows 2-7 down to 3-8
$ws.Rows.Item(2).Copy() | Out-Null
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new scan record
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "231249"
$ws.Range("B2").Value = "Anatomy"
$ws.Range("C2").Value = "14/08/2025"
$ws.Range("D2").Value = "10:38:02"
$ws.Range("E2").Value = "Scan"
$ws.Range("F2").Value = "admin@admin.com"

# Adjust the Log Time values that differ from a plain shift-down
$ws.Range("D3").Value = "10:38:02"
$ws.Range("D6").Value = "10:38:03"
